$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the two new SKU values (commit: "Add files via upload / 10006046 / 10102228")
$ws.Range("A106").Value = 10006046
$ws.Range("A107").Value = 10102228

# Normalize the font used throughout column A (incl. the new rows) to Calibri 10,
# matching the reformat captured in the diff (existing per-cell colors are preserved).
$ws.Range("A1:A107").Font.Name = "Calibri"
$ws.Range("A1:A107").Font.Size = 10

# Restore the view state: scrolled near the bottom of the list with G100 selected.
$ws.Range("A92").Select()
$ws.Range("G100").Select()
